$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3090
$ws1.Range("F3").Value = 497
$ws1.Range("F4").Value = 65
$ws1.Range("F5").Value = 56
$ws1.Range("F6").Value = 10
$ws1.Range("F8").Value = 17
$ws1.Range("F9").Value = 1075
$ws1.Range("F10").Value = 15077
$ws1.Range("F11").Value = 196
$ws1.Range("F12").Value = 148
$ws1.Range("F13").Value = 508
$ws1.Range("F14").Value = 5982
$ws1.Range("F17").Value = 54
$ws1.Range("F19").Value = 1249
$ws1.Range("F20").Value = 24
$ws1.Range("F22").Value = 5
$ws1.Range("F24").Value = 830
$ws1.Range("F25").Value = 2964
$ws1.Range("F27").Value = 10813
$ws1.Range("F28").Value = 1219
$ws1.Range("F29").Value = 92
$ws1.Range("F30").Value = 134
$ws1.Range("F31").Value = 3764
$ws1.Range("F32").Value = 254
$ws1.Range("F33").Value = 71

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3090
$ws4.Range("F4").Value = 497
$ws4.Range("F5").Value = 65
$ws4.Range("F6").Value = 56
$ws4.Range("F7").Value = 10
$ws4.Range("F9").Value = 17
$ws4.Range("F10").Value = 1075
$ws4.Range("F11").Value = 15077
$ws4.Range("F12").Value = 196
$ws4.Range("F13").Value = 148
$ws4.Range("F14").Value = 508
$ws4.Range("F15").Value = 5982
$ws4.Range("F18").Value = 54
$ws4.Range("F20").Value = 1249
$ws4.Range("F21").Value = 24
$ws4.Range("F23").Value = 5
$ws4.Range("F25").Value = 830
$ws4.Range("F26").Value = 2964
$ws4.Range("F29").Value = 10813
$ws4.Range("F30").Value = 1219
$ws4.Range("F31").Value = 92
$ws4.Range("F32").Value = 134
$ws4.Range("F33").Value = 3764
$ws4.Range("F34").Value = 254
$ws4.Range("F35").Value = 71
